$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 47-81: overwrite hword1 (A), hword2 (B) and hcategory (D) with the
# corrected/deduplicated dataset entries. hcorrAns (C) stays "left" for all
# of these rows and is left untouched.
$updates = @(
    @(47, "बहस", "blaokbonvd", 10),
    @(48, "कंबल", "poldcemps", 10),
    @(49, "नल", "ltbster", 10),
    @(50, "गाय", "valuabzr", 10),
    @(51, "गाजर", "dasewn", 10),
    @(52, "सावधान", "capmpiwity", 10),
    @(53, "गुफा", "perigrm", 10),
    @(54, "मछली", "ucter", 10),
    @(55, "फूल", "suggeyf", 10),
    @(56, "मगरम", "hodker", 8),
    @(57, "यूोप", "builde", 8),
    @(58, "गहा", "chulk", 8),
    @(59, "मदान", "gway", 8),
    @(60, "नाजु", "ceok", 8),
    @(61, "अनुबध", "leab", 8),
    @(62, "पतल", "dival", 8),
    @(63, "बाइन", "atnton", 8),
    @(64, "अजमदा", "uttac", 8),
    @(65, "वकी", "averake", 8),
    @(66, "केल", "fgor", 8),
    @(67, "नक", "bnat", 8),
    @(68, "असमल", "beaufifsl", 8),
    @(69, "संट", "subtraction", 9),
    @(70, "टराव", "tumor", 9),
    @(71, "परावान", "octopus", 9),
    @(72, "मगरम", "bow", 9),
    @(73, "यूोप", "mustache", 9),
    @(74, "अजमदा", "tricycle", 9),
    @(75, "केल", "cauliflower", 9),
    @(76, "असमल", "rug", 9),
    @(77, "अपरा", "blackboard", 9),
    @(78, "संट", "kitchen", 9),
    @(79, "जगल", "policeman", 9),
    @(80, "टराव", "lobster", 9),
    @(81, "परावान", "valuable", 9)
)

foreach ($row in $updates) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Rows 82-101 were duplicate/erroneous entries missed by mistake; delete
# them entirely so the table shrinks from A1:D101 down to A1:D81.
$ws.Range("A82:D101").EntireRow.Delete() | Out-Null

# Restore the view/selection state recorded for the sheet after the edit.
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("D67").Select() | Out-Null
